$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3048
$ws1.Range("F3").Value = 733
$ws1.Range("F5").Value = 6789
$ws1.Range("F6").Value = 1793
$ws1.Range("F7").Value = 29
$ws1.Range("F10").Value = 63
$ws1.Range("F11").Value = 131
$ws1.Range("F12").Value = 143

# Sheet "全部类型" (all types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3049
$ws4.Range("F4").Value = 733
$ws4.Range("F6").Value = 6789
$ws4.Range("F7").Value = 1793
$ws4.Range("F8").Value = 29
$ws4.Range("F11").Value = 63
$ws4.Range("F12").Value = 131
$ws4.Range("F13").Value = 143
